# Add a new TPA ("MEDIASSIST") block to the extraction-details table.
# The existing table is a repeating 5-column-per-TPA layout (APPROVAL /
# ENHANCEMENT / DISCHARGE under "APPROVAL", then DENIED, then QUERY), with
# row 1 = TPA name, row 2 = status category, row 3 = sub-header, rows 4-15 =
# field values. We copy the last existing block (STAR, columns AU:AY) into
# the five brand-new columns AZ:BD so the borders/fills/fonts/merges are
# carried over exactly, then overwrite the copied text with the new TPA's
# labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New columns' widths (AZ..BD = columns 52..56) -----------------
$ws.Columns.Item(52).ColumnWidth = 41.9
$ws.Columns.Item(53).ColumnWidth = 41.5
$ws.Columns.Item(54).ColumnWidth = 41.6
$ws.Columns.Item(55).ColumnWidth = 35.9
$ws.Columns.Item(56).ColumnWidth = 55.5

# --- 2. Clone formatting (borders/fills/fonts/merged cells) from the
#        last TPA block (STAR, AU:AY) onto the new block (AZ:BD) ------
$ws.Range("AU1:AY3").Copy($ws.Range("AZ1:BD3"))
$ws.Range("AU4:AY15").Copy($ws.Range("AZ4:BD15"))

# --- 3. Header text for the new "MEDIASSIST" TPA ----------------------
$ws.Range("AZ1").Value = "MEDIASSIST"
$ws.Range("AZ2").Value = "APPROVAL"
$ws.Range("BC2").Value = "DENIED"
$ws.Range("BD2").Value = "QUERY"
$ws.Range("AZ3").Value = "PRE AUTH "
$ws.Range("BA3").Value = "ENHANCEMENT"
$ws.Range("BB3").Value = "DISCHARGE"

# --- 4. Field values (rows 4-15) for the new TPA -----------------------
# Row 4 - ID
$ws.Range("AZ4").Value = "Medi Assist ID"
$ws.Range("BA4").Value = "Medi Assist ID"
$ws.Range("BB4").Value = "Medi Assist ID"
$ws.Range("BC4").Value = "Medi Assist ID"
$ws.Range("BD4").Value = "Medi Assist ID"

# Row 5 - CLAIM NUMBER
$ws.Range("AZ5").Value = "In brackets after `"Cashless Authorization Letter`""
$ws.Range("BA5").Value = "In brackets after `"Cashless Authorization Letter`""
$ws.Range("BB5").Value = "In brackets after `"Cashless Authorization Letter`""
$ws.Range("BC5").Value = "Cashless Claim Reference Number"
$ws.Range("BD5").Value = "Cashless Claim Reference Number"

# Row 6 - NAME
$ws.Range("AZ6").Value = "Patient Name"
$ws.Range("BA6").Value = "Patient Name"
$ws.Range("BB6").Value = "Patient Name"
$ws.Range("BC6").Value = "Patient Name"
$ws.Range("BD6").Value = "Patient Name"

# Row 7 - POLICY NO
$ws.Range("AZ7").Value = "Policy No."
$ws.Range("BA7").Value = "Policy No."
$ws.Range("BB7").Value = "Policy No."
$ws.Range("BC7").Value = "Policy No."
$ws.Range("BD7").Value = "Policy No."

# Row 8 - POLICY PERIOD
$ws.Range("AZ8").Value = "Policy Period"
$ws.Range("BA8").Value = "Policy/Plan Period"
$ws.Range("BB8").Value = "Policy/Plan Period"
$ws.Range("BC8").Value = "null"
$ws.Range("BD8").Value = "null"

# Row 9 - ROHINI ID
$ws.Range("AZ9").Value = "Rohini Id"
$ws.Range("BA9").Value = "Rohini Id"
$ws.Range("BB9").Value = "Rohini Id"
$ws.Range("BC9").Value = "null"
$ws.Range("BD9").Value = "null"

# Row 10 - DATE OF ADMISSION
$ws.Range("AZ10").Value = "Expected Date Of Admission"
$ws.Range("BA10").Value = "Expected Date Of Admission"
$ws.Range("BB10").Value = "Expected/Actual Date Of Admission"
$ws.Range("BC10").Value = "null"
$ws.Range("BD10").Value = "null"

# Row 11 - DATE OF DISCHARGE
$ws.Range("AZ11").Value = "Estimated Date of Discharge"
$ws.Range("BA11").Value = "Estimated/Actual Date of Discharge"
$ws.Range("BB11").Value = "Estimated/Actual Date of Discharge"
$ws.Range("BC11").Value = "null"
$ws.Range("BD11").Value = "null"

# Row 12 - AMOUNT
$ws.Range("AZ12").Value = "Total Authorized amount "
$ws.Range("BA12").Value = "Total Authorized amount "
$ws.Range("BB12").Value = "Total Authorized amount "
$ws.Range("BC12").Value = "null"
$ws.Range("BD12").Value = "null"

# Row 13 - NOTE
$ws.Range("AZ13").Value = "Authorization Remarks"
$ws.Range("BA13").Value = "Authorization Remarks"
$ws.Range("BB13").Value = "Authorization Remarks"
$ws.Range("BC13").Value = "Next line after `"We regret to inform you`""
$ws.Range("BD13").Value = "Next line after `" We require the following additional information`""

# Row 14 - HOSPITAL ADDRESS
$ws.Range("AZ14").Value = "Next line after `"To`""
$ws.Range("BA14").Value = "Next line after `"To`""
$ws.Range("BB14").Value = "Next line after `"To`""
$ws.Range("BC14").Value = "Next line after `"To`""
$ws.Range("BD14").Value = "Next line after `"To`""

# Row 15 - DATE AND TIME
$ws.Range("AZ15").Value = "Date :"
$ws.Range("BA15").Value = "Date :"
$ws.Range("BB15").Value = "Date :"
$ws.Range("BC15").Value = "Date :"
$ws.Range("BD15").Value = "Date :"

# --- 5. Small unrelated label fix in the existing STAR block ----------
$ws.Range("AX14").Value = "Next line after `"To`""

# --- 6. Move the view/selection over to the newly-added block ---------
$ws.Range("BB20").Select()
